$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(4)

# Header row (row 1): B..M
$headerArr = New-Object 'object[,]' 1,12
$headerArr[0,0] = "bank"
$headerArr[0,1] = "deposit_type"
$headerArr[0,2] = "currency"
$headerArr[0,3] = "owner"
$headerArr[0,4] = "total"
$headerArr[0,5] = "property_category"
$headerArr[0,6] = "category"
$headerArr[0,7] = "date"
$headerArr[0,8] = "legislator_name"
$headerArr[0,9] = "legislator_id"
$headerArr[0,10] = "source_file"
$headerArr[0,11] = "index"
$ws.Range("B1:M1").Value = $headerArr

# Data rows 2..12: columns B..M
$dataArr = New-Object 'object[,]' 11,12
$dataArr[0,0] = "臺灣銀行群賢分行"
$dataArr[0,1] = "活期儲蓄存款"
$dataArr[0,2] = "新臺幣"
$dataArr[0,3] = "管碧玲"
$dataArr[0,4] = 2181267
$dataArr[0,5] = "deposit"
$dataArr[0,6] = "normal"
$dataArr[0,7] = "2011-12-16"
$dataArr[0,8] = "管碧玲"
$dataArr[0,9] = 1374
$dataArr[0,10] = "tmp32301"
$dataArr[0,11] = 58
$dataArr[1,0] = "合作金庫商業銀行西門分行"
$dataArr[1,1] = "活期儲蓄存款"
$dataArr[1,2] = "新臺幣"
$dataArr[1,3] = "管碧玲"
$dataArr[1,4] = 109999
$dataArr[1,5] = "deposit"
$dataArr[1,6] = "normal"
$dataArr[1,7] = "2011-12-16"
$dataArr[1,8] = "管碧玲"
$dataArr[1,9] = 1374
$dataArr[1,10] = "tmp32301"
$dataArr[1,11] = 59
$dataArr[2,0] = "高雄銀行市府分行"
$dataArr[2,1] = "活期存款"
$dataArr[2,2] = "新臺幣"
$dataArr[2,3] = "管碧玲"
$dataArr[2,4] = 61239
$dataArr[2,5] = "deposit"
$dataArr[2,6] = "normal"
$dataArr[2,7] = "2011-12-16"
$dataArr[2,8] = "管碧玲"
$dataArr[2,9] = 1374
$dataArr[2,10] = "tmp32301"
$dataArr[2,11] = 60
$dataArr[3,0] = "國泰世華商業銀行中正分行"
$dataArr[3,1] = "活期儲蓄存款"
$dataArr[3,2] = "新臺幣"
$dataArr[3,3] = "管碧玲"
$dataArr[3,4] = 417453
$dataArr[3,5] = "deposit"
$dataArr[3,6] = "normal"
$dataArr[3,7] = "2011-12-16"
$dataArr[3,8] = "管碧玲"
$dataArr[3,9] = 1374
$dataArr[3,10] = "tmp32301"
$dataArr[3,11] = 61
$dataArr[4,0] = "高雄市府郵局(第1g支局)"
$dataArr[4,1] = "活期儲蓄存款"
$dataArr[4,2] = "新臺幣"
$dataArr[4,3] = "管碧玲"
$dataArr[4,4] = 238177
$dataArr[4,5] = "deposit"
$dataArr[4,6] = "normal"
$dataArr[4,7] = "2011-12-16"
$dataArr[4,8] = "管碧玲"
$dataArr[4,9] = 1374
$dataArr[4,10] = "tmp32301"
$dataArr[4,11] = 62
$dataArr[5,0] = "合作金庫商業銀行十全分行"
$dataArr[5,1] = "活期儲蓄存款"
$dataArr[5,2] = "新臺幣"
$dataArr[5,3] = "管碧玲"
$dataArr[5,4] = 994
$dataArr[5,5] = "deposit"
$dataArr[5,6] = "normal"
$dataArr[5,7] = "2011-12-16"
$dataArr[5,8] = "管碧玲"
$dataArr[5,9] = 1374
$dataArr[5,10] = "tmp32301"
$dataArr[5,11] = 63
$dataArr[6,0] = "台北富邦商業銀行北投分行"
$dataArr[6,1] = "活期儲蓄存款"
$dataArr[6,2] = "新臺幣"
$dataArr[6,3] = "管碧玲"
$dataArr[6,4] = 1238
$dataArr[6,5] = "deposit"
$dataArr[6,6] = "normal"
$dataArr[6,7] = "2011-12-16"
$dataArr[6,8] = "管碧玲"
$dataArr[6,9] = 1374
$dataArr[6,10] = "tmp32301"
$dataArr[6,11] = 65
$dataArr[7,0] = "陽信商業銀行大屯分行"
$dataArr[7,1] = "活期儲蓄存款"
$dataArr[7,2] = "新臺幣"
$dataArr[7,3] = "管碧玲"
$dataArr[7,4] = 9792
$dataArr[7,5] = "deposit"
$dataArr[7,6] = "normal"
$dataArr[7,7] = "2011-12-16"
$dataArr[7,8] = "管碧玲"
$dataArr[7,9] = 1374
$dataArr[7,10] = "tmp32301"
$dataArr[7,11] = 66
$dataArr[8,0] = "合作金庫商業銀行長春分行"
$dataArr[8,1] = "活期儲蓄存款"
$dataArr[8,2] = "新臺幣"
$dataArr[8,3] = "管碧玲"
$dataArr[8,4] = 30904
$dataArr[8,5] = "deposit"
$dataArr[8,6] = "normal"
$dataArr[8,7] = "2011-12-16"
$dataArr[8,8] = "管碧玲"
$dataArr[8,9] = 1374
$dataArr[8,10] = "tmp32301"
$dataArr[8,11] = 67
$dataArr[9,0] = "臺灣銀行北投分行"
$dataArr[9,1] = "綜合存款"
$dataArr[9,2] = "新臺幣"
$dataArr[9,3] = "管碧玲"
$dataArr[9,4] = 107480
$dataArr[9,5] = "deposit"
$dataArr[9,6] = "normal"
$dataArr[9,7] = "2011-12-16"
$dataArr[9,8] = "管碧玲"
$dataArr[9,9] = 1374
$dataArr[9,10] = "tmp32301"
$dataArr[9,11] = 68
$dataArr[10,0] = "合作金庫商業銀行營業部"
$dataArr[10,1] = "綜合存款"
$dataArr[10,2] = "美金"
$dataArr[10,3] = "管碧玲"
$dataArr[10,4] = 122250.13
$dataArr[10,5] = "deposit"
$dataArr[10,6] = "normal"
$dataArr[10,7] = "2011-12-16"
$dataArr[10,8] = "管碧玲"
$dataArr[10,9] = 1374
$dataArr[10,10] = "tmp32301"
$dataArr[10,11] = 69
$ws.Range("B2:M12").Value = $dataArr

# Apply header style (copy format from B1, which already has style s=1) to new header cells F1:M1
$ws.Range("B1").Copy()
$ws.Range("F1:M1").PasteSpecial(-4122)

# Apply data-row style (copy format from B2, which already has style s=2) to new data cells F2:M12, and to F (total) cells
$ws.Range("B2").Copy()
$ws.Range("F2:M12").PasteSpecial(-4122)

$excel.CutCopyMode = 0
